$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the M8 note: the old "Sequences not annotated by IMGT..." note is replaced
# by a new note about sequences not identified by Watson et al. This also removes
# the old shared string and appends the new text (as the last unique string).
$ws.Range("M8").Value = "Sequences notidentified by Watson et al. in this assembly but marked as ORF/P in other sources"

# TRAJ row (row 14) counts updated
$ws.Range("B14").Value = 50
$ws.Range("D14").Value = 5

# TRGJ row (row 22) counts updated, and the "non-canonical J-motif" note (M22) removed
$ws.Range("B22").Value = 5
$ws.Range("C22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()

# Move the active selection from M15 to D26
$ws.Range("D26").Select()
